$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B10").Value = "www.stat.gov.kg"
